$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 377
$ws1.Range("F4").Value = 197
$ws1.Range("F5").Value = 18
$ws1.Range("F6").Value = 1242
$ws1.Range("F7").Value = 449
$ws1.Range("F9").Value = 188
$ws1.Range("F11").Value = 173
$ws1.Range("F12").Value = 1052
$ws1.Range("F13").Value = 4
$ws1.Range("F15").Value = 187
$ws1.Range("F16").Value = 1495
$ws1.Range("F17").Value = 548
$ws1.Range("F18").Value = 229
$ws1.Range("F21").Value = 817
$ws1.Range("F22").Value = 1150
$ws1.Range("F24").Value = 1916
$ws1.Range("F25").Value = 2658
$ws1.Range("F26").Value = 1435
$ws1.Range("F27").Value = 62
$ws1.Range("F28").Value = 36
$ws1.Range("F29").Value = 408
$ws1.Range("F30").Value = 416
$ws1.Range("F31").Value = 1227
$ws1.Range("F32").Value = 821
$ws1.Range("F33").Value = 1354
$ws1.Range("F34").Value = 162
$ws1.Range("F36").Value = 784
$ws1.Range("F37").Value = 607
$ws1.Range("F38").Value = 671
$ws1.Range("F39").Value = 845
$ws1.Range("F40").Value = 363
$ws1.Range("F41").Value = 245

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 7
$ws2.Range("F15").Value = 626

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 377
$ws4.Range("F7").Value = 197
$ws4.Range("F8").Value = 18
$ws4.Range("F9").Value = 7
$ws4.Range("F11").Value = 1242
$ws4.Range("F12").Value = 449
$ws4.Range("F14").Value = 188
$ws4.Range("F17").Value = 173
$ws4.Range("F18").Value = 1052
$ws4.Range("F20").Value = 187
$ws4.Range("F21").Value = 1495
$ws4.Range("F22").Value = 548
$ws4.Range("F23").Value = 229
$ws4.Range("F28").Value = 1150
$ws4.Range("F29").Value = 2658
$ws4.Range("F30").Value = 1435
$ws4.Range("F31").Value = 62
$ws4.Range("F34").Value = 408
$ws4.Range("F35").Value = 416
$ws4.Range("F36").Value = 1227
$ws4.Range("F39").Value = 821
$ws4.Range("F40").Value = 1354
$ws4.Range("F41").Value = 784
$ws4.Range("F42").Value = 607
$ws4.Range("F43").Value = 671
$ws4.Range("F44").Value = 845
$ws4.Range("F45").Value = 363
$ws4.Range("F48").Value = 245
